$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top of the data, shifting everything down
$ws.Rows.Item(1).Insert()

# Set the new cell's value (shared string already exists elsewhere in the sheet)
$ws.Range("A1").Value = "GO_Biological_Process_2023"

# Move the active cell/selection to H20, matching the post-edit view state
$ws.Range("H20").Select()
